$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    $headers = $sec.Headers
    foreach ($hdr in $headers) {
        if ($hdr.Exists) {
            $hdr.Range.Find.Execute("December 2022", $true, $false, $false, $false, $false,
                                     $true, 1, $false, "August 2020", 2)
            foreach ($shp in $hdr.Shapes) {
                $tf = $shp.TextFrame
                if ($tf.HasText) {
                    $tf.TextRange.Find.Execute("December 2022", $true, $false, $false, $false, $false,
                                                $true, 1, $false, "August 2020", 2)
                }
            }
        }
    }
}
